# Update countries & provincias Spain
# Applies the data refresh captured in the commit "Update countries & provincias Spain":
#  - update counters for several countries (simple in-place updates)
#  - update Malta's counters which causes it to drop below Jordania / Georgia /
#    Republica de Chipre / Letonia in the (descending) sort by "Casos totales"
#  - swap the order of "Montserrat" and "Islas Malvinas" (same totals, order fix)
#  - bump the "Datos actualizados" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Simple in-place counter updates (country keeps its row/rank) ----

# Row 6 - India
$ws.Range("B6").Value = 2594612
$ws.Range("C6").Value = 5404
$ws.Range("D6").Value = 1863239
$ws.Range("E6").Value = 681247
$ws.Range("G6").Value = 42
$ws.Range("H6").Value = 50126

# Row 22 - Alemania
$ws.Range("B22").Value = 224562
$ws.Range("C22").Value = 84
$ws.Range("E22").Value = 12372

# Row 45 - Paises Bajos
$ws.Range("B45").Value = 63002
$ws.Range("C45").Value = 507
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 6172

# Row 62 - Azerbaiyan
$ws.Range("B62").Value = 34219
$ws.Range("C62").Value = 112
$ws.Range("D62").Value = 31875
$ws.Range("E62").Value = 1838

# Row 73 - El Salvador
$ws.Range("B73").Value = 22912
$ws.Range("C73").Value = 293
$ws.Range("D73").Value = 10807
$ws.Range("E73").Value = 11493

# Row 79 - Dinamarca
$ws.Range("B79").Value = 15617
$ws.Range("C79").Value = 134
$ws.Range("D79").Value = 13340
$ws.Range("E79").Value = 1656

# Row 103 - Croacia
$ws.Range("B103").Value = 6571
$ws.Range("C103").Value = 151
$ws.Range("D103").Value = 5220
$ws.Range("E103").Value = 1185
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = 166

# Row 124 - Sri Lanka
$ws.Range("B124").Value = 2893
$ws.Range("C124").Value = 3
$ws.Range("D124").Value = 2670
$ws.Range("E124").Value = 212

# Row 151 - Burkina Faso
$ws.Range("B151").Value = 1249
$ws.Range("C151").Value = 9
$ws.Range("D151").Value = 1013
$ws.Range("E151").Value = 182

# ---- Malta block: Malta's update drops it below Jordania, Georgia, ----
# ---- Republica de Chipre and Letonia, so rows 144-148 get re-sorted ----
# New (sorted) order for rows 144-148, keeping each country's own data
# except Malta which receives fresh updated counters.

# Row 144 -> Jordania (was row 145's data)
$ws.Range("A144").Value = "Jordania"
$ws.Range("B144").Value = 1339
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 1229
$ws.Range("E144").Value = 99
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 11

# Row 145 -> Georgia (was row 146's data)
$ws.Range("A145").Value = "Georgia"
$ws.Range("B145").Value = 1336
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 1088
$ws.Range("E145").Value = 231
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 17

# Row 146 -> Republica de Chipre (was row 147's data)
$ws.Range("A146").Value = "Republica de Chipre"
$ws.Range("B146").Value = 1332
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 870
$ws.Range("E146").Value = 442
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 20

# Row 147 -> Letonia (was row 148's data)
$ws.Range("A147").Value = "Letonia"
$ws.Range("B147").Value = 1322
$ws.Range("C147").Value = 7
$ws.Range("D147").Value = 1078
$ws.Range("E147").Value = 212
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 32

# Row 148 -> Malta (now with its own freshly updated counters)
$ws.Range("A148").Value = "Malta"
$ws.Range("B148").Value = 1306
$ws.Range("C148").Value = 63
$ws.Range("D148").Value = 749
$ws.Range("E148").Value = 548
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 9

# ---- Swap "Islas Malvinas" and "Montserrat" (same total, order fix) ----

# Row 213 -> Montserrat (was "Islas Malvinas")
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214 -> Islas Malvinas (was "Montserrat")
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# ---- Bump the "last updated" timestamp shown in A1 ----
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 14:52"
